# Weekly update for "Fruta, Terminal Hortofrutícola Agro Chillán - Mandarina":
# insert two new daily observations (Murcott, Primera & Segunda, 2021-10-25 i.e.
# serial 44494, sold "$/caja 18 kilos") at the top of the data block (row 68),
# pushing the existing rows 68:159 down to 70:161.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 68 (shifts 68..159 -> 70..161).
$ws.Rows.Item(68).Insert()
$ws.Rows.Item(68).Insert()

# New row 68: Murcott / Primera
$ws.Cells.Item(68, 1).Value  = 7
$ws.Cells.Item(68, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(68, 3).Value  = "Ñuble"
$ws.Cells.Item(68, 4).Value  = 44494
$ws.Cells.Item(68, 5).Value  = 16
$ws.Cells.Item(68, 6).Value  = "Fruta"
$ws.Cells.Item(68, 7).Value  = 100102
$ws.Cells.Item(68, 8).Value  = "Cítricos"
$ws.Cells.Item(68, 9).Value  = 100102004
$ws.Cells.Item(68, 10).Value = "Mandarina"
$ws.Cells.Item(68, 11).Value = "Murcott"
$ws.Cells.Item(68, 12).Value = "Primera"
$ws.Cells.Item(68, 13).Value = 120
$ws.Cells.Item(68, 14).Value = 7000
$ws.Cells.Item(68, 15).Value = 7500
$ws.Cells.Item(68, 16).Value = 7250
$ws.Cells.Item(68, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(68, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(68, 19).Value = 403
$ws.Cells.Item(68, 20).Value = 18

# New row 69: Murcott / Segunda
$ws.Cells.Item(69, 1).Value  = 7
$ws.Cells.Item(69, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(69, 3).Value  = "Ñuble"
$ws.Cells.Item(69, 4).Value  = 44494
$ws.Cells.Item(69, 5).Value  = 16
$ws.Cells.Item(69, 6).Value  = "Fruta"
$ws.Cells.Item(69, 7).Value  = 100102
$ws.Cells.Item(69, 8).Value  = "Cítricos"
$ws.Cells.Item(69, 9).Value  = 100102004
$ws.Cells.Item(69, 10).Value = "Mandarina"
$ws.Cells.Item(69, 11).Value = "Murcott"
$ws.Cells.Item(69, 12).Value = "Segunda"
$ws.Cells.Item(69, 13).Value = 120
$ws.Cells.Item(69, 14).Value = 6000
$ws.Cells.Item(69, 15).Value = 6500
$ws.Cells.Item(69, 16).Value = 6250
$ws.Cells.Item(69, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(69, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(69, 19).Value = 347
$ws.Cells.Item(69, 20).Value = 18
